$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A: sequential ID numbers for rows 4..17 (continuing 1,2 already in A2:A3)
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10
$ws.Range("A12").Value = 11
$ws.Range("A13").Value = 12
$ws.Range("A14").Value = 13
$ws.Range("A15").Value = 14
$ws.Range("A16").Value = 15
$ws.Range("A17").Value = 16

# Columns W (Patente), X (Motor), Y (Chasis): renumber test data MJK066..MJK079
$ws.Range("W4").Value = "MJK066"
$ws.Range("X4").Value = "ABCD0RRGA066"
$ws.Range("Y4").Value = "ZXC0987RGA066"

$ws.Range("W5").Value = "MJK067"
$ws.Range("X5").Value = "ABCD0RRGA067"
$ws.Range("Y5").Value = "ZXC0987RGA067"

$ws.Range("W6").Value = "MJK068"
$ws.Range("X6").Value = "ABCD0RRGA068"
$ws.Range("Y6").Value = "ZXC0987RGA068"

$ws.Range("W7").Value = "MJK069"
$ws.Range("X7").Value = "ABCD0RRGA069"
$ws.Range("Y7").Value = "ZXC0987RGA069"

$ws.Range("W8").Value = "MJK070"
$ws.Range("X8").Value = "ABCD0RRGA070"
$ws.Range("Y8").Value = "ZXC0987RGA070"

$ws.Range("W9").Value = "MJK071"
$ws.Range("X9").Value = "ABCD0RRGA071"
$ws.Range("Y9").Value = "ZXC0987RGA071"

$ws.Range("W10").Value = "MJK072"
$ws.Range("X10").Value = "ABCD0RRGA072"
$ws.Range("Y10").Value = "ZXC0987RGA072"

$ws.Range("W11").Value = "MJK073"
$ws.Range("X11").Value = "ABCD0RRGA073"
$ws.Range("Y11").Value = "ZXC0987RGA073"

$ws.Range("W12").Value = "MJK074"
$ws.Range("X12").Value = "ABCD0RRGA074"
$ws.Range("Y12").Value = "ZXC0987RGA074"

$ws.Range("W13").Value = "MJK075"
$ws.Range("X13").Value = "ABCD0RRGA075"
$ws.Range("Y13").Value = "ZXC0987RGA075"

$ws.Range("W14").Value = "MJK076"
$ws.Range("X14").Value = "ABCD0RRGA076"
$ws.Range("Y14").Value = "ZXC0987RGA076"

$ws.Range("W15").Value = "MJK077"
$ws.Range("X15").Value = "ABCD0RRGA077"
$ws.Range("Y15").Value = "ZXC0987RGA077"

$ws.Range("W16").Value = "MJK078"
$ws.Range("X16").Value = "ABCD0RRGA078"
$ws.Range("Y16").Value = "ZXC0987RGA078"

$ws.Range("W17").Value = "MJK079"
$ws.Range("X17").Value = "ABCD0RRGA079"
$ws.Range("Y17").Value = "ZXC0987RGA079"

# Reset the view: select the whole of row 4 (also clears the saved top-left scroll cell)
$ws.Rows.Item(4).Select()
